$d = $word.ActiveDocument

# Replace the text of a whole paragraph (matched by its trimmed text),
# preserving any run-level attributes Word would keep on an in-place
# edit (e.g. xml:space="preserve"), and preserving any comment anchors
# (commentRangeStart/commentRangeEnd/commentReference) that sit in the
# paragraph alongside the text run.
function Replace-ParaText($oldText, $newText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($ptext -eq $oldText) {
            $p.Range.Text = $newText
            return $true
        }
    }
    return $false
}

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# "English" occurs twice: once inside the language-switcher hyperlink
# (must stay "English") and once as the standalone heading that names
# this translation (must become "Englisch"). A Range scoped to just the
# standalone paragraph doesn't persist reliably, so replace both
# occurrences and then restore the hyperlink's text.
Replace-Text "English" "Englisch"
Replace-Text "English" "Englisch"
$hyperlink = $d.Hyperlinks.Item(1)
$hyperlink.Range.Text = "English"

Replace-ParaText "To confirm your registration, we would require you and one guest of your choice to provide us with:" "Um Ihre Anmeldung zu bestätigen, bitten wir Sie und einen Gast Ihrer Wahl, uns folgende Angaben zu machen:" | Out-Null
Replace-ParaText "A scanned copy of your international passports" "Eine gescannte Kopie Ihrer internationalen Pässe" | Out-Null

# This run is immediately followed by a comment range end / reference;
# use Find so the replacement run ends up without xml:space="preserve".
Replace-Text "Covid-19 vaccination certificates" "Covid-19-Impfbescheinigungen"

Replace-ParaText "Your country manager will be in touch to confirm your booking or request any other relevant details. " "Ihr Ländermanager wird sich mit Ihnen in Verbindung setzen, um Ihre Buchung zu bestätigen oder weitere Details zu erfragen. " | Out-Null
Replace-ParaText "Our event package offers you and your guest: " "Unser Veranstaltungspaket bietet Ihnen und Ihren Gästen: " | Out-Null
Replace-ParaText "Flight tickets " "Flugtickets " | Out-Null
Replace-ParaText "Travel insurance " "Reiseversicherung " | Out-Null
Replace-ParaText "Airport – Hotel – Airport transfer " "Flughafen - Hotel - Flughafentransfer " | Out-Null
Replace-ParaText "One hotel room for you and your guest / Two hotel rooms for you and your guest" "Ein Hotelzimmer für Sie und Ihren Gast / Zwei Hotelzimmer für Sie und Ihren Gast" | Out-Null
Replace-ParaText "Meals (Breakfast, lunch, and dinner)" "Mahlzeiten (Frühstück, Mittag- und Abendessen)" | Out-Null
Replace-ParaText "We will send you a confirmation letter before your departure date with the event agenda and information about your flights, transportation, and accommodation. " "Wir senden Ihnen vor Ihrer Abreise ein Bestätigungsschreiben mit dem Veranstaltungsprogramm und Informationen zu Ihren Flügen, Ihrem Transport und Ihrer Unterkunft. " | Out-Null
Replace-ParaText "We look forward to seeing you soon." "Wir freuen uns darauf, Sie bald wiederzusehen." | Out-Null

Write-Output "Done"
